$d = $word.ActiveDocument

# Locate the paragraph that contains the astromap credit/link line:
# " Jenika Hollana, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
# by scanning paragraphs for the unique "Jenika" marker, and confirm it
# still references the 2018 map before touching it.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text
    if (($t -like "*Jenika*") -and ($t -like "*GaNight/2018*")) {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $pEnd = $target.Range.End

    # Remove all existing run content in the paragraph (but keep the
    # paragraph mark / pPr intact) by deleting everything up to, but not
    # including, the trailing paragraph mark.
    $body = $d.Range($pStart, $pEnd - 1)
    $body.Delete()

    # Insert the replacement text (year updated 2018 -> 2022) as a single
    # new run with no special character formatting, collapsing all of the
    # old per-word runs/spell-check markers into one plain run.
    $newText = " Jenika Hollana, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
    $insertRange = $d.Range($pStart, $pStart)
    $insertRange.InsertAfter($newText)
}
